# Powerpoint writer: avoid extra blank lines before author.
# (In the case where there is no subtitle.)
# Remove the now-unused "Subtitle 2" placeholder shape from slide 1 — it
# only contained two empty line breaks, which is what the writer used to
# emit when no subtitle text was actually present.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The "Subtitle 2" shape is a layout-backed placeholder (subTitle, idx=1),
# so the first Delete() only resets it back to an empty, layout-default
# placeholder instance (re-seated/renamed, e.g. "Subtitle 3") rather than
# removing it from the slide's shape tree. Deleting it again (same
# position, shape #2) removes it for good, leaving just the title shape -
# matching a slide that never had a subtitle.
$s.Shapes.Item(2).Delete()
$s.Shapes.Item(2).Delete()
